# "Generate Report for Handoff" — refresh the localization-status report:
#   * Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#     on the Overview sheet (per-language columns) and on each language sheet.
#   * The "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" stamps
#     advance to the new generation time.
#   * The (now shorter) status text no longer needs as wide a column, so the
#     status columns are narrowed.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = $newStatus   # Overview: zh-cn status
$wsOverview.Range("F2").Value = $newStatus   # Overview: de-de status
$wsZhCn.Range("C2").Value     = $newStatus   # zh-cn: Status
$wsDeDe.Range("C2").Value     = $newStatus   # de-de: Status

# --- Handoff timestamps bump forward to the new generation run ---
$wsOverview.Range("G2").Value = "2016-09-07 11:23:50"   # Latest HO Xliff Generate Date
$wsZhCn.Range("H2").Value     = "2016-09-07 11:23:45"   # zh-cn: Latest Handoff Datetime
$wsDeDe.Range("H2").Value     = "2016-09-07 11:23:50"   # de-de: Latest Handoff Datetime

# --- Narrow the status columns (they held the long "Handed back..." text) ---
# Stored column width = ColumnWidth + 5/6, so subtract 5/6 to land on the
# target stored width of 17.2159881591797 characters.
$newColumnWidth = 17.2159881591797 - (5 / 6)

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth   # Overview column E
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth   # Overview column F
$wsZhCn.Columns.Item(3).ColumnWidth     = $newColumnWidth   # zh-cn column C
$wsDeDe.Columns.Item(3).ColumnWidth     = $newColumnWidth   # de-de column C
